$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-01-12 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-01-13 Saturday", 2)

# Update the division problems inside the table. Each data row (1, 5, 9, 13, 17)
# of the single table holds 5 problems, one per column; the other rows are blank
# answer rows. The new values for several cells coincide with the original (or
# new) values of other cells, so instead of a document-wide Find/Replace (which
# can match the wrong occurrence once duplicates exist) each cell's Range.Text
# is set directly -- this only ever touches that specific cell.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "66÷9=" },
    @{ Row = 1;  Col = 2; New = "73÷7=" },
    @{ Row = 1;  Col = 3; New = "21÷5=" },
    @{ Row = 1;  Col = 4; New = "26÷6=" },
    @{ Row = 1;  Col = 5; New = "55÷5=" },

    @{ Row = 5;  Col = 1; New = "89÷7=" },
    @{ Row = 5;  Col = 2; New = "10÷5=" },
    @{ Row = 5;  Col = 3; New = "42÷5=" },
    @{ Row = 5;  Col = 4; New = "27÷8=" },
    @{ Row = 5;  Col = 5; New = "66÷9=" },

    @{ Row = 9;  Col = 1; New = "56÷7=" },
    @{ Row = 9;  Col = 2; New = "66÷2=" },
    @{ Row = 9;  Col = 3; New = "30÷7=" },
    @{ Row = 9;  Col = 4; New = "89÷4=" },
    @{ Row = 9;  Col = 5; New = "76÷9=" },

    @{ Row = 13; Col = 1; New = "97÷3=" },
    @{ Row = 13; Col = 2; New = "24÷9=" },
    @{ Row = 13; Col = 3; New = "85÷2=" },
    @{ Row = 13; Col = 4; New = "46÷6=" },
    @{ Row = 13; Col = 5; New = "84÷3=" },

    @{ Row = 17; Col = 1; New = "83÷2=" },
    @{ Row = 17; Col = 2; New = "39÷9=" },
    @{ Row = 17; Col = 3; New = "96÷3=" },
    @{ Row = 17; Col = 4; New = "87÷5=" },
    @{ Row = 17; Col = 5; New = "47÷6=" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
